# "Generate Report for Handoff"
# A new handoff cycle happened for e2e\b.md: the file is now "Ready for
# handoff" in both locales, the handoff xliff files were regenerated
# (new content hash, new timestamp), the old handback file/version is no
# longer the latest (Content Duplicate -> False) and an explanatory
# Error Detail message is now populated. The Error Detail column is also
# widened so the message is readable.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# ---- Overview sheet: row 3 is the e2e\b.md entry ----
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-17 10:35:11"

# ---- zh-cn sheet: row 3 is the b.md entry ----
$wsZhCn.Range("C3").Value = "Ready for handoff"
# Leading "'" forces Excel to store "False" as text rather than a boolean;
# the Style reset afterwards clears the resulting quote-prefix formatting
# so the cell keeps the plain/default style of its neighbours.
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-17 10:35:00"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/67378b97842fc3df53fc5395ed1f1d9127203966/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/611c810593fabebd3d47122226e297a5dc196648/e2e/b.md."
# ColumnWidth goes through a char<->pixel rounding conversion, so 39.14
# (not 40) is what round-trips to a stored column width of exactly 40.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.14

# ---- de-de sheet: row 3 is the b.md entry ----
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-17 10:35:11"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/67378b97842fc3df53fc5395ed1f1d9127203966/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/611c810593fabebd3d47122226e297a5dc196648/e2e/b.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.14
